$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns C (runs), D (balls), E (fours), F (sixes)
# for rows 2..11, as per the updated activity data.
$data = @{
    2  = @("19", "20", "4", "0")
    3  = @("10", "7",  "1", "1")
    4  = @("23", "24", "2", "0")
    5  = @("36", "28", "7", "0")
    6  = @("16", "19", "0", "1")
    7  = @("0",  "3",  "0", "0")
    8  = @("25", "15", "2", "2")
    9  = @("53", "48", "2", "1")
    10 = @("97", "55", "7", "6")
    11 = @("61", "43", "6", "2")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($col = 3; $col -le 6; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $vals[$col - 3]
    }
}
